$wb = $excel.ActiveWorkbook
$admin = $wb.Worksheets.Item("Admin")
$pim = $wb.Worksheets.Item("PIM")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Copy formatting first (does not touch the shared-string table) ---

# Row 3 headers - copy formatting from existing header cells
$admin.Range("A1").Copy()
$admin.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$admin.Range("B1").Copy()
$admin.Range("B3:D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$admin.Range("C1").Copy()
$admin.Range("E3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$admin.Range("E1").Copy()
$admin.Range("F3:G3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$dashboard.Range("F2").Copy()
$admin.Range("F4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$admin.Range("A1").Copy()
$admin.Range("H4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$admin.Range("B4").ClearFormats()

# --- Now set the cell values, in authoring order ---

# rename existing test case
$admin.Range("A1").Value = "TC01_Admin_SearchUserByNameAndRole"
$admin.Range("F2").Value = "TC01_Admin_SearchUserByNameAndRole"

# new test case header + data
$admin.Range("A3").Value = "TC02_Admin_CreateSystemUser"
$admin.Range("F3").Value = "Status"
$admin.Range("F4").Value = "Enabled"
$admin.Range("G3").Value = "UserPassword"
$admin.Range("B3").Value = "EmployeeFirstName"
$admin.Range("C3").Value = "EmployeeMiddleName"
$admin.Range("D3").Value = "EmployeeLastName"
$admin.Range("B4").Value = "Priya"
$admin.Range("C4").Value = "Raj"
$admin.Range("D4").Value = "Sharma"
$admin.Range("E4").Value = "Admin"
$admin.Range("E3").Value = "UserRole"
$admin.Range("G4").Value = "R29vZHdpbGwxMjM0NQ=="
$admin.Range("H4").Value = "TC02_Admin_CreateSystemUser"

# --- PIM sheet: C10 loses its stray fill-only format ---
$pim.Range("C10").ClearFormats()
$pim.Range("C10").Value = "Sharma"

# --- Selections: update PIM's selection without stealing the active tab ---
$pim.Range("C5:C6").Select()
$admin.Activate()
$admin.Range("N4").Select()
